$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("TestProgram_Counter")

# ------------------------------------------------------------------
# Start from a clean slate: wipe all existing cell content/formatting
# and column-level formatting so nothing bleeds into the new layout.
# ------------------------------------------------------------------
$ws3.Cells.Clear()
$ws3.Columns.Item(2).ClearFormats()
$ws3.Columns.Item(3).ClearFormats()

# Rename the sheet
$ws3.Name = "Test Programs"

# ------------------------------------------------------------------
# Column widths (character units as accepted by ColumnWidth); the
# engine snaps to the nearest renderable pixel width.
# ------------------------------------------------------------------
$ws3.Columns.Item(3).ColumnWidth = 18.1    # -> stored width 19
$ws3.Columns.Item(4).ColumnWidth = 16.92   # -> stored width ~17.71 (bestFit)
$ws3.Columns.Item(5).ColumnWidth = 16.92   # -> stored width ~17.71 (bestFit)

# ------------------------------------------------------------------
# Cell numeric / formula content (does not affect the shared-string
# table, so it can be filled in any order).
# ------------------------------------------------------------------
$ws3.Range("A3").Value = 0
$ws3.Range("A4").Value = 3
$ws3.Range("A5").Value = 6
$ws3.Range("A6").Value = 9
$ws3.Range("B3:B6").Formula = "=DEC2HEX(A3,3)"

$ws3.Range("A9").Value = 6
$ws3.Range("B9").Formula = "=DEC2HEX(A9,3)"

$ws3.Range("A17").Value = 0
$ws3.Range("A18").Value = 3
$ws3.Range("A19").Value = 6
$ws3.Range("A20").Value = 9
$ws3.Range("A21").Value = 12
$ws3.Range("A22").Value = 15
$ws3.Range("B17:B22").Formula = "=DEC2HEX(A17,3)"

$ws3.Range("A23").Value = 18
$ws3.Range("B23").Formula = "=DEC2HEX(A23,3)"

$ws3.Range("A24").Value = 21
$ws3.Range("A25").Value = 24
$ws3.Range("A26").Value = 27
$ws3.Range("B24:B26").Formula = "=DEC2HEX(A24,3)"

# ------------------------------------------------------------------
# Text content, entered in the same order as the original authoring
# session so the shared-string table indices line up.
# ------------------------------------------------------------------
$ws3.Range("C17").Value = "LD B, 0x1"
$ws3.Range("C18").Value = "LD E, 0x0"
$ws3.Range("C19").Value = "LD F, 0xff"
$ws3.Range("C20").Value = "LD A, E"
$ws3.Range("C21").Value = "ADD A, B"
$ws3.Range("C22").Value = "LD E, A"
$ws3.Range("C23").Value = "LD A, F"
$ws3.Range("C24").Value = "SUB A, B"
$ws3.Range("C25").Value = "LD F, A"
$ws3.Range("C26").Value = "JP [0x009]"

$ws3.Range("C3").Value = "LD B, 0x1"
$ws3.Range("C4").Value = "LD A, 0x20"
$ws3.Range("C5").Value = "ADD A, B"
$ws3.Range("C6").Value = "JP [0x006]"
$ws3.Range("C9").Value = "SUB A, B"

$ws3.Range("E20").Value = "Equiv to ADD E, B"
$ws3.Range("E20:E22").NumberFormat = "@"
$ws3.Range("E20:E22").VerticalAlignment = -4108
$ws3.Range("E20:E22").HorizontalAlignment = -4108
$ws3.Range("E20:E22").Merge()

$ws3.Range("D17").Value = "04 80 01"
$ws3.Range("D18").Value = "06 00 00"
$ws3.Range("D19").Value = "06 80 ff"
$ws3.Range("D20").Value = "08 40 00"
$ws3.Range("D21").Value = "14 10 00"
$ws3.Range("D23").Value = "08 50 00"
$ws3.Range("D25").Value = "0a 80 00"
$ws3.Range("D22").Value = "0a 00 00"
$ws3.Range("D26").Value = "2c 00 09"
$ws3.Range("D17:D26").NumberFormat = "@"

$ws3.Range("D3").Value = "04 80 01"
$ws3.Range("D4").Value = "04 00 20"
$ws3.Range("D5").Value = "14 10 00"
$ws3.Range("D6").Value = "2c 00 06"
$ws3.Range("D9").Value = "18 10 00"
$ws3.Range("D3:D6").NumberFormat = "@"
$ws3.Range("D9").NumberFormat = "@"

$ws3.Range("D24").Value = "18 10 00"
$ws3.Range("D24").NumberFormat = "@"

$ws3.Range("C16").Value = "// Increments E, decrements F"
$ws3.Range("C16").Font.ThemeColor = 5

$ws3.Range("E23").Value = "Equiv to SUB F, B"
$ws3.Range("E23:E25").NumberFormat = "@"
$ws3.Range("E23:E25").VerticalAlignment = -4108
$ws3.Range("E23:E25").HorizontalAlignment = -4108
$ws3.Range("E23:E25").Merge()

$ws3.Range("C2").Value = "// Increments A"
$ws3.Range("C2").Font.ThemeColor = 5

$ws3.Range("A1").Value = "Dec"
$ws3.Range("B1").Value = "Hexa"
$ws3.Range("A1:B1").Font.Bold = $true
$ws3.Range("A2:B2").Font.Bold = $true

# ------------------------------------------------------------------
# View: selection over the first merged comment cell
# ------------------------------------------------------------------
$ws3.Range("E20:E22").Select()
